$wb = $excel.ActiveWorkbook

# --- Sheet 1: Predictions NO Tuning ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("C4").Value = 0.4986999471354288
$ws1.Range("D4").Value = 0.5013000528645712
$ws1.Range("B5").Value = "U.N.C."
$ws1.Range("B6").Value = "U.N.C."
$ws1.Range("C7").Value = 0.5213918732611434
$ws1.Range("D7").Value = 0.4786081267388564
$ws1.Range("B8").Value = "U.N.C."
$ws1.Range("B9").Value = "U.N.C."
$ws1.Range("B10").Value = "U.N.C."
$ws1.Range("C11").Value = 0.5692839823501323
$ws1.Range("D11").Value = 0.4307160176498677
$ws1.Range("B16").Value = "U.N.C."
$ws1.Range("C18").Value = 0.5310794522270156
$ws1.Range("D18").Value = 0.4689205477729845
$ws1.Range("C22").Value = 0.5447250124667524
$ws1.Range("D22").Value = 0.4552749875332476
$ws1.Range("C23").Value = 0.537870978919868
$ws1.Range("D23").Value = 0.4621290210801319
$ws1.Range("B24").Value = "U.N.C."
$ws1.Range("B25").Value = "U.N.C."
$ws1.Range("B26").Value = "U.N.C."
$ws1.Range("C28").Value = 0.491102133507783
$ws1.Range("D28").Value = 0.5088978664922169
$ws1.Range("B31").Value = "U.N.C."
$ws1.Range("C33").Value = 0.4718533505002138
$ws1.Range("D33").Value = 0.5281466494997862
$ws1.Range("B34").Value = "U.N.C."
$ws1.Range("B36").Value = "U.N.C."
$ws1.Range("C37").Value = 0.5411531160569628
$ws1.Range("D37").Value = 0.4588468839430372
$ws1.Range("B38").Value = "U.N.C."
$ws1.Range("C39").Value = 0.4605997705688807
$ws1.Range("D39").Value = 0.5394002294311193
$ws1.Range("C41").Value = 0.5817175396086411
$ws1.Range("D41").Value = 0.4182824603913588
$ws1.Range("C42").Value = 0.5281978769981704
$ws1.Range("D42").Value = 0.4718021230018296

# --- Sheet 2: Predictions WITH Tuning ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("C4").Value = 0.4523889073091529
$ws2.Range("D4").Value = 0.5476110926908471
$ws2.Range("B5").Value = "U.N.C."
$ws2.Range("B6").Value = "U.N.C."
$ws2.Range("C7").Value = 0.4860478629936003
$ws2.Range("D7").Value = 0.5139521370063997
$ws2.Range("B8").Value = "U.N.C."
$ws2.Range("B9").Value = "U.N.C."
$ws2.Range("B10").Value = "U.N.C."
$ws2.Range("C11").Value = 0.5400631779254161
$ws2.Range("D11").Value = 0.4599368220745839
$ws2.Range("B16").Value = "U.N.C."
$ws2.Range("C18").Value = 0.498049458061175
$ws2.Range("D18").Value = 0.5019505419388249
$ws2.Range("C22").Value = 0.514760354218787
$ws2.Range("D22").Value = 0.485239645781213
$ws2.Range("C23").Value = 0.5153052426198728
$ws2.Range("D23").Value = 0.4846947573801272
$ws2.Range("B24").Value = "U.N.C."
$ws2.Range("B25").Value = "U.N.C."
$ws2.Range("B26").Value = "U.N.C."
$ws2.Range("C28").Value = 0.4477647865174925
$ws2.Range("D28").Value = 0.5522352134825075
$ws2.Range("B31").Value = "U.N.C."
$ws2.Range("C33").Value = 0.4334053548128298
$ws2.Range("D33").Value = 0.5665946451871702
$ws2.Range("B34").Value = "U.N.C."
$ws2.Range("B36").Value = "U.N.C."
$ws2.Range("C37").Value = 0.5058599895644404
$ws2.Range("D37").Value = 0.4941400104355596
$ws2.Range("B38").Value = "U.N.C."
$ws2.Range("C39").Value = 0.4189023727609251
$ws2.Range("D39").Value = 0.581097627239075
$ws2.Range("C41").Value = 0.5589542030282891
$ws2.Range("D41").Value = 0.4410457969717108
$ws2.Range("C42").Value = 0.4897837774194581
$ws2.Range("D42").Value = 0.5102162225805418

# --- Sheet 3: Safe Seats ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B4").Value = "U.N.C."
$ws3.Range("B5").Value = "U.N.C."
$ws3.Range("B6").Value = "U.N.C."
$ws3.Range("B7").Value = "U.N.C."
$ws3.Range("B8").Value = "U.N.C."
$ws3.Range("B13").Value = "U.N.C."
$ws3.Range("B18").Value = "U.N.C."
$ws3.Range("B19").Value = "U.N.C."
$ws3.Range("B20").Value = "U.N.C."
$ws3.Range("B24").Value = "U.N.C."
$ws3.Range("B26").Value = "U.N.C."
$ws3.Range("B28").Value = "U.N.C."
$ws3.Range("B29").Value = "U.N.C."

# --- Sheet 4: Marginal Seats NO Tuning ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("C2").Value = 0.4986999471354288
$ws4.Range("D2").Value = 0.5013000528645712
$ws4.Range("C3").Value = 0.5213918732611434
$ws4.Range("D3").Value = 0.4786081267388564
$ws4.Range("C4").Value = 0.5692839823501323
$ws4.Range("D4").Value = 0.4307160176498677
$ws4.Range("C5").Value = 0.5310794522270156
$ws4.Range("D5").Value = 0.4689205477729845
$ws4.Range("C6").Value = 0.5447250124667524
$ws4.Range("D6").Value = 0.4552749875332476
$ws4.Range("C7").Value = 0.537870978919868
$ws4.Range("D7").Value = 0.4621290210801319
$ws4.Range("C8").Value = 0.491102133507783
$ws4.Range("D8").Value = 0.5088978664922169
$ws4.Range("C9").Value = 0.4718533505002138
$ws4.Range("D9").Value = 0.5281466494997862
$ws4.Range("C10").Value = 0.5411531160569628
$ws4.Range("D10").Value = 0.4588468839430372
$ws4.Range("C11").Value = 0.4605997705688807
$ws4.Range("D11").Value = 0.5394002294311193
$ws4.Range("C12").Value = 0.5817175396086411
$ws4.Range("D12").Value = 0.4182824603913588
$ws4.Range("C13").Value = 0.5281978769981704
$ws4.Range("D13").Value = 0.4718021230018296

# --- Sheet 5: Marginal Seats Tuning ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("C2").Value = 0.4523889073091529
$ws5.Range("D2").Value = 0.5476110926908471
$ws5.Range("C3").Value = 0.4860478629936003
$ws5.Range("D3").Value = 0.5139521370063997
$ws5.Range("C4").Value = 0.5400631779254161
$ws5.Range("D4").Value = 0.4599368220745839
$ws5.Range("C5").Value = 0.498049458061175
$ws5.Range("D5").Value = 0.5019505419388249
$ws5.Range("C6").Value = 0.514760354218787
$ws5.Range("D6").Value = 0.485239645781213
$ws5.Range("C7").Value = 0.5153052426198728
$ws5.Range("D7").Value = 0.4846947573801272
$ws5.Range("C8").Value = 0.4477647865174925
$ws5.Range("D8").Value = 0.5522352134825075
$ws5.Range("C9").Value = 0.4334053548128298
$ws5.Range("D9").Value = 0.5665946451871702
$ws5.Range("C10").Value = 0.5058599895644404
$ws5.Range("D10").Value = 0.4941400104355596
$ws5.Range("C11").Value = 0.4189023727609251
$ws5.Range("D11").Value = 0.581097627239075
$ws5.Range("C12").Value = 0.5589542030282891
$ws5.Range("D12").Value = 0.4410457969717108
$ws5.Range("C13").Value = 0.4897837774194581
$ws5.Range("D13").Value = 0.5102162225805418
